$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 15 (phone "09876543" stored as text, birthday "2020-08-16",
# points 0) needs to be split into two rows:
#   - row 15: the phone recorded as a plain number (9876543), same birthday, points 0
#   - row 16: the original text phone "09876543" (leading zero kept), no birthday, points 0
#
# Stash the current row 15 in a scratch row first so we can copy its exact cell
# values/types back out afterwards (Copy preserves the original text/number typing
# instead of letting Excel re-interpret strings such as dates or leading-zero numbers).
$ws.Range("A15:C15").Copy($ws.Range("A30:C30"))

# Build new row 16 from the stashed values: keep the original text phone number
# and points, but leave the birthday blank (B10 is a known blank cell, used purely
# as a source so the destination cell stays a genuine empty cell instead of being
# removed, which is what happens when assigning an empty string directly).
$ws.Range("A30").Copy($ws.Cells.Item(16, 1))
$ws.Range("B10").Copy($ws.Cells.Item(16, 2))
$ws.Range("C30").Copy($ws.Cells.Item(16, 3))

# Build new row 15: phone becomes a real number, birthday is carried over unchanged,
# points stay at 0.
$ws.Cells.Item(15, 1).Value2 = 9876543
$ws.Range("B30").Copy($ws.Cells.Item(15, 2))
$ws.Cells.Item(15, 3).Value2 = 0

# Clean up the scratch row used to stage the original values.
$ws.Range("A30:C30").ClearContents()
